# Apply the diet-code cleanup edit:
# 1) Fix a handful of inconsistent IBMR codes (the "20 different versions of the diet code")
# 2) Sort the data range by PostCategory (column A) ascending
# 3) Re-point view/selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Step 1: fix inconsistent / duplicate IBMR code spellings -------------
# Crangon spp: "shrimp" -> "SHRIMP"
$ws.Range("B25").Value = "SHRIMP"
# Exopalaemon modestus: "shrimp" -> "other"
$ws.Range("B41").Value = "other"
# Sinocalanus spp: "sino" -> "othcalad"
$ws.Range("B79").Value = "othcalad"

# --- Step 2: sort A2:C108 by column A ascending ----------------------------
$rng = $ws.Range("A2:C108")
$key1 = $ws.Range("A2:A108")
$rng.Sort($key1, 1, $null, $null, 1, $null, 1, 1, $false, $null, $null, 1)

# --- Step 3: update sheet view ---------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("C1").Select()
